# Reorder each year's 12 monthly rows so that Oct/Nov/Dec come first,
# followed by Jan..Sep, for every year block in the sheet.
#
# Data layout: row 1 = headers, rows 2-49 = 4 years (2014-2017) x 12 months,
# columns A (date label) through J (9 numeric series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
$numCols = 10
$yearBlockSize = 12
$numYearBlocks = 4

$blockStart = $startRow
for ($b = 0; $b -lt $numYearBlocks; $b++) {
    # Read this year's 12 rows (all columns) into a jagged array.
    $rows = @()
    for ($r = $blockStart; $r -lt ($blockStart + $yearBlockSize); $r++) {
        $row = @()
        for ($c = 1; $c -le $numCols; $c++) {
            $row += $ws.Cells.Item($r, $c).Value()
        }
        $rows += ,$row
    }

    # Rotate: months at index 9,10,11 (Oct,Nov,Dec) move to the front,
    # followed by months at index 0..8 (Jan..Sep).
    $rotated = @()
    for ($i = 9; $i -lt 12; $i++) { $rotated += ,$rows[$i] }
    for ($i = 0; $i -lt 9; $i++) { $rotated += ,$rows[$i] }

    # Write the rotated rows back in place.
    for ($i = 0; $i -lt $yearBlockSize; $i++) {
        $r = $blockStart + $i
        $rowData = $rotated[$i]
        for ($c = 1; $c -le $numCols; $c++) {
            $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
        }
    }

    $blockStart += $yearBlockSize
}
